$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new poll data rows (70-82): Ifop, OpinionWay, and Ipsos polls (update 10/22)
# Row 70
$ws.Range("A70").Value = 21
$ws.Range("B70").Value = 2021
$ws.Range("C70").Value = 7
$ws.Range("D70").Value = 10
$ws.Range("E70").Value = 15
$ws.Range("F70").Value = "ifop"
$ws.Range("G70").Value = "online"
$ws.Range("H70").Value = "included"
$ws.Range("I70").Value = 1182
$ws.Range("J70").Value = 0.5
$ws.Range("K70").Value = 0.5
$ws.Range("L70").Value = 8
$ws.Range("M70").Value = 2
$ws.Range("N70").Value = 1.5
$ws.Range("O70").Value = 7
$ws.Range("P70").Value = 6
$ws.Range("Q70").Value = 24
$ws.Range("T70").Value = 15
$ws.Range("U70").Value = 0.5
$ws.Range("V70").Value = 2
$ws.Range("W70").Value = 17
$ws.Range("X70").Value = 16

# Row 71
$ws.Range("A71").Value = 21
$ws.Range("B71").Value = 2021
$ws.Range("C71").Value = 7
$ws.Range("D71").Value = 10
$ws.Range("E71").Value = 15
$ws.Range("F71").Value = "ifop"
$ws.Range("G71").Value = "online"
$ws.Range("H71").Value = "included"
$ws.Range("I71").Value = 1182
$ws.Range("J71").Value = 0.5
$ws.Range("K71").Value = 0.5
$ws.Range("L71").Value = 8
$ws.Range("M71").Value = 2
$ws.Range("N71").Value = 2
$ws.Range("O71").Value = 7
$ws.Range("P71").Value = 6
$ws.Range("Q71").Value = 25
$ws.Range("S71").Value = 10
$ws.Range("U71").Value = 1
$ws.Range("V71").Value = 3
$ws.Range("W71").Value = 18
$ws.Range("X71").Value = 17

# Row 72
$ws.Range("A72").Value = 21
$ws.Range("B72").Value = 2021
$ws.Range("C72").Value = 7
$ws.Range("D72").Value = 10
$ws.Range("E72").Value = 15
$ws.Range("F72").Value = "ifop"
$ws.Range("G72").Value = "online"
$ws.Range("H72").Value = "included"
$ws.Range("I72").Value = 1182
$ws.Range("J72").Value = 1
$ws.Range("K72").Value = 0.5
$ws.Range("L72").Value = 8.5
$ws.Range("M72").Value = 2.5
$ws.Range("N72").Value = 2
$ws.Range("O72").Value = 7
$ws.Range("P72").Value = 6.5
$ws.Range("Q72").Value = 24
$ws.Range("R72").Value = 10
$ws.Range("U72").Value = 1
$ws.Range("V72").Value = 3
$ws.Range("W72").Value = 17
$ws.Range("X72").Value = 17

# Row 73
$ws.Range("A73").Value = 22
$ws.Range("B73").Value = 2021
$ws.Range("C73").Value = 8
$ws.Range("D73").Value = 10
$ws.Range("E73").Value = 19
$ws.Range("F73").Value = "opinionway"
$ws.Range("G73").Value = "online"
$ws.Range("H73").Value = "included"
$ws.Range("I73").Value = 859
$ws.Range("J73").Value = 2
$ws.Range("K73").Value = 1
$ws.Range("L73").Value = 9
$ws.Range("M73").Value = 2
$ws.Range("N73").Value = 2
$ws.Range("O73").Value = 8
$ws.Range("P73").Value = 5
$ws.Range("Q73").Value = 26
$ws.Range("T73").Value = 9
$ws.Range("V73").Value = 4
$ws.Range("W73").Value = 19
$ws.Range("X73").Value = 13

# Row 74
$ws.Range("A74").Value = 22
$ws.Range("B74").Value = 2021
$ws.Range("C74").Value = 8
$ws.Range("D74").Value = 10
$ws.Range("E74").Value = 19
$ws.Range("F74").Value = "opinionway"
$ws.Range("G74").Value = "online"
$ws.Range("H74").Value = "included"
$ws.Range("I74").Value = 838
$ws.Range("J74").Value = 2
$ws.Range("K74").Value = "T_1"
$ws.Range("L74").Value = 10
$ws.Range("M74").Value = 2
$ws.Range("N74").Value = 3
$ws.Range("O74").Value = 9
$ws.Range("P74").Value = 5
$ws.Range("Q74").Value = 26
$ws.Range("R74").Value = 6
$ws.Range("V74").Value = 4
$ws.Range("W74").Value = 19
$ws.Range("X74").Value = 14

# Row 75
$ws.Range("A75").Value = 22
$ws.Range("B75").Value = 2021
$ws.Range("C75").Value = 8
$ws.Range("D75").Value = 10
$ws.Range("E75").Value = 19
$ws.Range("F75").Value = "opinionway"
$ws.Range("G75").Value = "online"
$ws.Range("H75").Value = "included"
$ws.Range("I75").Value = 859
$ws.Range("J75").Value = 2
$ws.Range("K75").Value = 1
$ws.Range("L75").Value = 9
$ws.Range("M75").Value = 2
$ws.Range("N75").Value = 2
$ws.Range("O75").Value = 8
$ws.Range("P75").Value = 6
$ws.Range("Q75").Value = 26
$ws.Range("S75").Value = 7
$ws.Range("V75").Value = 4
$ws.Range("W75").Value = 20
$ws.Range("X75").Value = 14

# Row 76
$ws.Range("A76").Value = 22
$ws.Range("B76").Value = 2021
$ws.Range("C76").Value = 8
$ws.Range("D76").Value = 10
$ws.Range("E76").Value = 19
$ws.Range("F76").Value = "opinionway"
$ws.Range("G76").Value = "online"
$ws.Range("H76").Value = "included"
$ws.Range("I76").Value = 838
$ws.Range("J76").Value = 2
$ws.Range("K76").Value = "T_1"
$ws.Range("L76").Value = 10
$ws.Range("M76").Value = 2
$ws.Range("N76").Value = 3
$ws.Range("O76").Value = 8
$ws.Range("P76").Value = 5
$ws.Range("Q76").Value = 26
$ws.Range("T76").Value = 10
$ws.Range("V76").Value = 6
$ws.Range("W76").Value = 28

# Row 77
$ws.Range("A77").Value = 23
$ws.Range("B77").Value = 2021
$ws.Range("C77").Value = 6
$ws.Range("D77").Value = 10
$ws.Range("E77").Value = 10
$ws.Range("F77").Value = "ipsos"
$ws.Range("G77").Value = "online"
$ws.Range("H77").Value = "excluded"
$ws.Range("I77").Value = 8888
$ws.Range("J77").Value = 1
$ws.Range("K77").Value = 1
$ws.Range("L77").Value = 8
$ws.Range("M77").Value = 2
$ws.Range("N77").Value = 2
$ws.Range("O77").Value = 9
$ws.Range("P77").Value = 5
$ws.Range("Q77").Value = 24
$ws.Range("T77").Value = 13
$ws.Range("U77").Value = 1
$ws.Range("V77").Value = 3
$ws.Range("W77").Value = 15
$ws.Range("X77").Value = 16

# Row 78
$ws.Range("A78").Value = 23
$ws.Range("B78").Value = 2021
$ws.Range("C78").Value = 6
$ws.Range("D78").Value = 10
$ws.Range("E78").Value = 10
$ws.Range("F78").Value = "ipsos"
$ws.Range("G78").Value = "online"
$ws.Range("H78").Value = "excluded"
$ws.Range("I78").Value = 8796
$ws.Range("J78").Value = 1
$ws.Range("K78").Value = 1
$ws.Range("L78").Value = 8
$ws.Range("M78").Value = 2
$ws.Range("N78").Value = 2
$ws.Range("O78").Value = 9.5
$ws.Range("P78").Value = 5
$ws.Range("Q78").Value = 25.5
$ws.Range("R78").Value = 10
$ws.Range("U78").Value = 1
$ws.Range("V78").Value = 3
$ws.Range("W78").Value = 16
$ws.Range("X78").Value = 16

# Row 79
$ws.Range("A79").Value = 23
$ws.Range("B79").Value = 2021
$ws.Range("C79").Value = 6
$ws.Range("D79").Value = 10
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = "ipsos"
$ws.Range("G79").Value = "online"
$ws.Range("H79").Value = "excluded"
$ws.Range("I79").Value = 8714
$ws.Range("J79").Value = 1
$ws.Range("K79").Value = 1
$ws.Range("L79").Value = 8
$ws.Range("M79").Value = 2
$ws.Range("N79").Value = 2
$ws.Range("O79").Value = 9.5
$ws.Range("P79").Value = 5
$ws.Range("Q79").Value = 26
$ws.Range("S79").Value = 9
$ws.Range("U79").Value = 1
$ws.Range("V79").Value = 3
$ws.Range("W79").Value = 16
$ws.Range("X79").Value = 16.5

# Row 80
$ws.Range("A80").Value = 23
$ws.Range("B80").Value = 2021
$ws.Range("C80").Value = 6
$ws.Range("D80").Value = 10
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = "ipsos"
$ws.Range("G80").Value = "online"
$ws.Range("H80").Value = "excluded"
$ws.Range("I80").Value = 8698
$ws.Range("J80").Value = 1
$ws.Range("K80").Value = 1
$ws.Range("L80").Value = 8
$ws.Range("M80").Value = 2
$ws.Range("N80").Value = 2
$ws.Range("O80").Value = 9.5
$ws.Range("P80").Value = 5
$ws.Range("Q80").Value = 26
$ws.Range("T80").Value = 16
$ws.Range("U80").Value = 1
$ws.Range("V80").Value = 4.5
$ws.Range("W80").Value = 24

# Row 81
$ws.Range("A81").Value = 23
$ws.Range("B81").Value = 2021
$ws.Range("C81").Value = 6
$ws.Range("D81").Value = 10
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = "ipsos"
$ws.Range("G81").Value = "online"
$ws.Range("H81").Value = "excluded"
$ws.Range("I81").Value = 8541
$ws.Range("J81").Value = 1
$ws.Range("K81").Value = 1
$ws.Range("L81").Value = 8
$ws.Range("M81").Value = 2
$ws.Range("N81").Value = 2
$ws.Range("O81").Value = 10
$ws.Range("P81").Value = 5
$ws.Range("Q81").Value = 27
$ws.Range("R81").Value = 13
$ws.Range("U81").Value = 1
$ws.Range("V81").Value = 5
$ws.Range("W81").Value = 25

# Row 82
$ws.Range("A82").Value = 23
$ws.Range("B82").Value = 2021
$ws.Range("C82").Value = 6
$ws.Range("D82").Value = 10
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = "ipsos"
$ws.Range("G82").Value = "online"
$ws.Range("H82").Value = "excluded"
$ws.Range("I82").Value = 8541
$ws.Range("J82").Value = 1
$ws.Range("K82").Value = 1
$ws.Range("L82").Value = 8
$ws.Range("M82").Value = 2
$ws.Range("N82").Value = 2
$ws.Range("O82").Value = 10
$ws.Range("P82").Value = 5
$ws.Range("Q82").Value = 28
$ws.Range("S82").Value = 12
$ws.Range("U82").Value = 1
$ws.Range("V82").Value = 5
$ws.Range("W82").Value = 25


# Update the sheet view: freeze header row, scroll to show the new rows,
# and set the active selection to match the saved view state.
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A67").Select()
$excel.ActiveWindow.ScrollRow = 67
$ws.Range("T71").Select()
